# Adapt bulac to fleet changes to Uruguay 2025 model
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FACTORS")

$values = @(
    @(0.1273426742845834, 0.1375864719616154, 0.009218662286240091),
    @(0.1385664810060629, 0.07361601420459957, 0.1608122037454987),
    @(0.02866171107854141, 0.1286148203061389, 0.1410991310337403),
    @(0.005599852723031887, 0.0439363938074464, 0.1681830836057347),
    @(0.08988188957996833, 0.1059079198305923, 0.02792260181143307),
    @(0.1108322445897598, 0.1721514704028065, 0.04982390814496991),
    @(0.07138001983692742, 0.0007266658150443441, 0.08715786115037481),
    @(0.04370309305202161, 0.1677989630732037, 0.07612722892127428),
    @(0.1618585534500533, 0.09055126348683162, 0.04144314431805703),
    @(0.049325297685688, 0.0522258868704585, 0.1307142713125733),
    @(0.1728481827133619, 0.02688413024126295, 0.107497903670104)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
    $ws.Cells.Item($row, 4).Value = $values[$i][2]
}
